$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2) values, matching the header columns in row 1:
# A=identifier  C=title  E=levelOfDescription  F=extentAndMedium  G=notes
$ws.Range("A2").Value = "MCH159-1"
$ws.Range("C2").Value = "NETHERLANDS REFORMED CHURCH, LETTERS & ATTACHMENTS,  BOOKS, DOCUMENTS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 22A | GRAP COUNT NUMER: NONE"

# D2 and H2 stay empty but are part of the formatted row (B2 is skipped;
# it keeps the sheet's default formatting).
$rng = $ws.Range("A2,C2:H2")
$rng.Font.ThemeColor = 1
$rng.Font.Name = "Calibri"
$rng.Font.Size = 10
